$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I62").Value = 5399.25
$ws.Range("N62").Value = -13059.625
$ws.Range("M62").Value = -4775.25
$ws.Range("J62").Value = 11811.625
$ws.Range("H62").Value = 9674.166999999999
$ws.Range("K62").Value = 5399.25
$ws.Range("L62").Value = 11811.625
$ws.Range("K65").Value = 26996.25
$ws.Range("H65").Value = 9674.166999999999
$ws.Range("N65").Value = -65298.125
$ws.Range("M65").Value = -23876.25
$ws.Range("I65").Value = 5399.25
$ws.Range("J65").Value = 11811.625
$ws.Range("L65").Value = 59058.125
$ws.Range("K137").Value = 40733.667
$ws.Range("H137").Value = 8316.556
$ws.Range("M137").Value = -38183.667
$ws.Range("I137").Value = 13577.889
$ws.Range("L138").Value = 11155.7799
$ws.Range("H138").Value = 3647.6143
$ws.Range("N138").Value = -21435.7799
$ws.Range("J138").Value = 3718.5933

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M45").Value = -1072.6
$ws.Range("K45").Value = 1449.6
$ws.Range("I45").Value = 1449.6
$ws.Range("H45").Value = 2282.2354
$ws.Range("I61").Value = 3193.5117
$ws.Range("H61").Value = 3096.6545
$ws.Range("M61").Value = -2981.5117
$ws.Range("K61").Value = 3193.5117
$ws.Range("I136").Value = 3193.5117
$ws.Range("K136").Value = 9580.535100000001
$ws.Range("H136").Value = 3096.6545
$ws.Range("M136").Value = -7030.535100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K5").Value = 1469.8
$ws.Range("L5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("M5").Value = -1356.8
$ws.Range("I5").Value = 1469.8
$ws.Range("N5").ClearContents()
$ws.Range("H5").Value = 1469.8
$ws.Range("I20").Value = 1005
$ws.Range("J20").Value = 1751.3334
$ws.Range("L20").Value = 1751.3334
$ws.Range("H20").Value = 1502.5555
$ws.Range("K20").Value = 1005
$ws.Range("M20").Value = -758
$ws.Range("N20").Value = -2245.3334
$ws.Range("J47").Value = 199999
$ws.Range("H47").Value = 199999
$ws.Range("N47").Value = -201039
$ws.Range("L47").Value = 199999
$ws.Range("N104").Value = -41988
$ws.Range("J104").Value = 35000
$ws.Range("H104").Value = 35000
$ws.Range("L104").Value = 35000
$ws.Range("M105").Value = 220.8571999999999
$ws.Range("I105").Value = 1526.1428
$ws.Range("K105").Value = 1526.1428
$ws.Range("H105").Value = 1743.7727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3320.5
$ws.Range("I31").Value = 1424.2354
$ws.Range("M31").Value = -1129.2354
$ws.Range("K31").Value = 1424.2354
$ws.Range("I34").Value = 1424.2354
$ws.Range("H34").Value = 3320.5
$ws.Range("M34").Value = -1222.2354
$ws.Range("K34").Value = 1424.2354
$ws.Range("M58").Value = -2171.8
$ws.Range("I58").Value = 2374.8
$ws.Range("K58").Value = 2374.8
$ws.Range("H58").Value = 2478.8333
$ws.Range("I99").Value = 2943.4
$ws.Range("H99").Value = 3318.2856
$ws.Range("M99").Value = -1445.4
$ws.Range("K99").Value = 2943.4
$ws.Range("I126").Value = 2943.4
$ws.Range("H126").Value = 3318.2856
$ws.Range("M126").Value = -6360.200000000001
$ws.Range("K126").Value = 8830.200000000001
$ws.Range("I136").Value = 2374.8
$ws.Range("K136").Value = 7124.400000000001
$ws.Range("H136").Value = 2478.8333
$ws.Range("M136").Value = -4574.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L39").Value = 36000
$ws.Range("N39").Value = -36588
$ws.Range("H39").Value = 6450
$ws.Range("J39").Value = 12000
$ws.Range("J86").Value = 370
$ws.Range("H86").Value = 446.53845
$ws.Range("L86").Value = 1110
$ws.Range("N86").Value = -3482
$ws.Range("I86").Value = 460.45456
$ws.Range("K86").Value = 1381.36368
$ws.Range("M86").Value = -195.3636799999999
$ws.Range("M89").Value = 1783.90896
$ws.Range("I89").Value = 460.45456
$ws.Range("L89").Value = 3330
$ws.Range("K89").Value = 4144.09104
$ws.Range("N89").Value = -15186
$ws.Range("H89").Value = 446.53845
$ws.Range("J89").Value = 370
$ws.Range("N92").Value = -5793.75
$ws.Range("M92").Value = -5764.125
$ws.Range("L92").Value = 3297.75
$ws.Range("J92").Value = 1099.25
$ws.Range("H92").Value = 1924.6666
$ws.Range("K92").Value = 7012.125
$ws.Range("I92").Value = 2337.375
$ws.Range("M113").Value = -2924
$ws.Range("H113").Value = 1280.7646
$ws.Range("K113").Value = 5094
$ws.Range("I113").Value = 1698
$ws.Range("I121").Value = 143503.86
$ws.Range("J121").Value = 884.1429000000001
$ws.Range("L121").Value = 2652.4287
$ws.Range("N121").Value = -5272.4287
$ws.Range("H121").Value = 72194
$ws.Range("M121").Value = -429201.58
$ws.Range("K121").Value = 430511.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I7").Value = 5000000
$ws.Range("K7").Value = 5000000
$ws.Range("M7").Value = -4999888
$ws.Range("H7").Value = 3100000
$ws.Range("H8").Value = 3100000
$ws.Range("I8").Value = 5000000
$ws.Range("K8").Value = 5000000
$ws.Range("M8").Value = -4999861
$ws.Range("J70").Value = 9017.5625
$ws.Range("L70").Value = 9017.5625
$ws.Range("H70").Value = 9268322
$ws.Range("N70").Value = -9557.5625
$ws.Range("J73").Value = 9017.5625
$ws.Range("H73").Value = 9268322
$ws.Range("L73").Value = 9017.5625
$ws.Range("N73").Value = -10889.5625
$ws.Range("H102").Value = 1313.75
$ws.Range("M102").Value = 659.8
$ws.Range("K102").Value = 962.2
$ws.Range("I102").Value = 962.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H5").Value = 0
$ws.Range("M40").Value = -7693.7856
$ws.Range("J40").Value = 8399.546
$ws.Range("N40").Value = -8671.546
$ws.Range("L40").Value = 8399.546
$ws.Range("K40").Value = 7829.7856
$ws.Range("H40").Value = 8080.48
$ws.Range("I40").Value = 7829.7856
$ws.Range("H74").Value = 47084.715
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61996
$ws.Range("L77").Value = 180000
$ws.Range("H77").Value = 47084.715
$ws.Range("N77").Value = -189984
$ws.Range("J77").Value = 60000
$ws.Range("I82").Value = 454477.6
$ws.Range("J82").Value = 4699.75
$ws.Range("M82").Value = -454116.6
$ws.Range("K82").Value = 454477.6
$ws.Range("H82").Value = 298033.12
$ws.Range("N82").Value = -5421.75
$ws.Range("L82").Value = 4699.75
$ws.Range("J85").Value = 4699.75
$ws.Range("N85").Value = -7195.75
$ws.Range("L85").Value = 4699.75
$ws.Range("H85").Value = 298033.12
$ws.Range("I85").Value = 454477.6
$ws.Range("M85").Value = -453229.6
$ws.Range("K85").Value = 454477.6
$ws.Range("N132").Value = -21410
$ws.Range("M132").Value = -15828.2
$ws.Range("K132").Value = 18358.2
$ws.Range("I132").Value = 6119.4
$ws.Range("H132").Value = 5928.143
$ws.Range("J132").Value = 5450
$ws.Range("L132").Value = 16350
$ws.Range("N136").Value = -29091
$ws.Range("I136").Value = 5749.6665
$ws.Range("K136").Value = 17248.9995
$ws.Range("H136").Value = 6070.7144
$ws.Range("J136").Value = 7997
$ws.Range("L136").Value = 23991
$ws.Range("M136").Value = -14698.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 967466.4399999999
$ws.Range("K2").Value = 967466.4399999999
$ws.Range("I2").Value = 967466.4399999999
$ws.Range("M2").Value = -967354.4399999999
$ws.Range("K11").Value = 50
$ws.Range("H11").Value = 50
$ws.Range("M11").Value = 92
$ws.Range("I11").Value = 50
$ws.Range("N14").Value = -21435.666
$ws.Range("K14").Value = 601.6
$ws.Range("J14").Value = 21099.666
$ws.Range("I14").Value = 601.6
$ws.Range("H14").Value = 8288.375
$ws.Range("L14").Value = 21099.666
$ws.Range("M14").Value = -433.6
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H107").Value = 1154
$ws.Range("M107").Value = -2145.273
$ws.Range("I107").Value = 1355.091
$ws.Range("K107").Value = 4065.273
$ws.Range("N109").Value = -143424
$ws.Range("H109").Value = 140650
$ws.Range("L109").Value = 140650
$ws.Range("J109").Value = 140650
$ws.Range("M113").Value = -1282.6001
$ws.Range("J113").Value = 2151.25
$ws.Range("H113").Value = 1498.826
$ws.Range("K113").Value = 3452.6001
$ws.Range("I113").Value = 1150.8667
$ws.Range("N113").Value = -10793.75
$ws.Range("L113").Value = 6453.75
$ws.Range("N120").Value = -180243.14
$ws.Range("H120").Value = 170567.14
$ws.Range("J120").Value = 170567.14
$ws.Range("L120").Value = 170567.14
$ws.Range("K122").Value = 13456.5
$ws.Range("M122").Value = -11006.5
$ws.Range("H122").Value = 6993.4443
$ws.Range("I122").Value = 4485.5
$ws.Range("M132").Value = -16029.5
$ws.Range("K132").Value = 18559.5
$ws.Range("I132").Value = 6186.5
$ws.Range("H132").Value = 5193.5
$ws.Range("I136").Value = 3108.9
$ws.Range("K136").Value = 9326.700000000001
$ws.Range("H136").Value = 3184.9285
$ws.Range("M136").Value = -6776.700000000001
